$d = $word.ActiveDocument

# --- Split the run containing "{m" into two runs: "{" and "m" ---
# (TokenIteratorFieldRewriterSplit breaks the field-open-delimiter "{" away
# from the following token text, each becoming its own run.)
$hit1 = $d.Content
$hit1.Find.Execute("{m", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($hit1.Find.Found) {
    $openBrace = $d.Range($hit1.Start, $hit1.Start + 1)
    # Force Word to split the run at this boundary without altering the
    # visible formatting: toggle Bold on then back off on the sub-range.
    $openBrace.Bold = 1
    $openBrace.Bold = 0
}

# --- Split the run containing ")}" into two runs: ")" and "}" ---
$hit2 = $d.Content
$hit2.Find.Execute(")}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($hit2.Find.Found) {
    # First break ")" away from "}" into its own run, keeping formatting.
    $closeParen = $d.Range($hit2.Start, $hit2.Start + 1)
    $closeParen.Bold = 1
    $closeParen.Bold = 0

    # Then rebuild the trailing "}" as a brand new run with no run
    # properties at all (matching the rewriter's output for the
    # field-close delimiter), by deleting it and re-inserting it fresh.
    $braceStart = $hit2.Start + 1
    $braceEnd = $hit2.Start + 2
    $closeBrace = $d.Range($braceStart, $braceEnd)
    $closeBrace.Text = ""
    $insertionPoint = $d.Range($braceStart, $braceStart)
    $insertionPoint.InsertAfter("}")
}
